$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.755.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.090.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.55%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.89'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.93%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.385'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0838'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.398.69'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.796'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.079.31'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.719.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0837'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.138'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.33%  '
$ws.Range("E29").Value = '  +11.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.17'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0614'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.12%  '
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("E38").Value = '  +2.52%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.28'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.539.25'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0223'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.57%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0926'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.29%  '
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.60%  '
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.10'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("E49").Value = '  +2.50%  '
$ws.Range("E50").Value = '  +1.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.286.86'
$ws.Range("D51").Style = "Normal"
